$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously empty "Total" results block (rows 30-32 / category row 17) ---
$ws.Range("C30").Value = 11
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = 50
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 0.56

$ws.Range("C31").Value = 86
$ws.Range("D31").Value = 88
$ws.Range("E31").Value = 841
$ws.Range("F31").Value = 9
$ws.Range("G31").Value = 45.1

$ws.Range("C32").Value = 318
$ws.Range("D32").Value = 320
$ws.Range("E32").Value = 2934
$ws.Range("F32").Value = 12
$ws.Range("G32").Value = 228.64

# --- Give the newly filled "F" column cells the same highlight style used by the
#     other category blocks (e.g. rows 15-17) ---
$ws.Range("F15").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F17").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Move the active selection to G32 (last cell touched) ---
$ws.Range("G32").Select() | Out-Null
